# map, stage 수치 조정 (map, stage numeric values adjustment)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Stage 4 ("map_size") row: tighten the map bounding box numbers.
$ws.Range("E4").Value = 145   # mapXmax 150 -> 145
$ws.Range("G4").Value = 130   # mapZmax 140 -> 130
$ws.Range("H4").Value = 72    # mapCenterX 75 -> 72
$ws.Range("I4").Value = 65    # mapCenterZ 70 -> 65

# Move the live selection from I11 to J11, matching the saved view state.
$ws.Range("J11").Select()
